# Apply crypto price/volume updates per the commit diff (Fri Jan 19 11:24:34 UTC 2024 refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.517.26"
$ws.Range("E2").Value = "  -2.37%  "
$ws.Range("D3").Value = "2.486.95"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'313.99"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "'94.79"
$ws.Range("E6").Value = "  -4.12%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("D10").Value = "'33.66"
$ws.Range("E10").Value = "  -4.46%  "
$ws.Range("D11").Value = "'0.0782"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'7.00"
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("D14").Value = "2.870.97"
$ws.Range("E14").Value = "  -1.36%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'15.50"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.481.24"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'0.795"
$ws.Range("E17").Value = "  -1.58%  "
$ws.Range("D18").Value = "41.481.42"
$ws.Range("E18").Value = "  -2.45%  "
$ws.Range("E19").Value = "  -3.65%  "
$ws.Range("D20").Value = "0.0₃0925"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'11.31"
$ws.Range("E21").Value = "  -6.62%  "
$ws.Range("D22").Value = "'68.98"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").Value = "'237.56"
$ws.Range("E23").Value = "  -1.33%  "
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'24.22"
$ws.Range("E27").Value = "  -4.37%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D30").Value = "'36.76"
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").Value = "'152.67"
$ws.Range("E31").Value = "  -2.44%  "
$ws.Range("D32").Value = "'5.51"
$ws.Range("E32").Value = "  -5.29%  "
$ws.Range("E33").Value = "  -3.58%  "
$ws.Range("D34").Value = "'18.17"
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("E35").Value = "  -3.59%  "
$ws.Range("E36").Value = "  -7.53%  "
$ws.Range("D37").Value = "'3.09"
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("D38").Value = "'1.88"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("E39").Value = "  -1.84%  "
$ws.Range("E40").Value = "  -6.82%  "
$ws.Range("D41").Value = "'4.24"
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "'19.92"
$ws.Range("E43").Value = "  -9.03%  "
$ws.Range("D44").Value = "1.999.49"
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").Value = "'3.04"
$ws.Range("E46").Value = "  -7.21%  "
$ws.Range("E47").Value = "  -2.57%  "
$ws.Range("D48").Value = "2.733.32"
$ws.Range("E48").Value = "  -1.00%  "
$ws.Range("E49").Value = "  -2.26%  "
$ws.Range("D50").Value = "'97.15"
$ws.Range("E50").Value = "  -3.45%  "
$ws.Range("E51").Value = "  -5.49%  "
